# Fruta / hortaliza, semanal
# Insert a new weekly record at row 87 (Macroferia Regional de Talca -
# Arándano (blue)), pushing the existing rows 87-98 down to 88-99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 87..98 down one position to make room for the new record.
$ws.Rows("87:87").Insert()

# Populate the newly inserted row 87 with the new weekly observation.
$ws.Cells.Item(87, 1).Value  = 5
$ws.Cells.Item(87, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(87, 3).Value  = "Maule"
$ws.Cells.Item(87, 4).Value  = 44918
$ws.Cells.Item(87, 5).Value  = 7
$ws.Cells.Item(87, 6).Value  = "Fruta"
$ws.Cells.Item(87, 7).Value  = 100101
$ws.Cells.Item(87, 8).Value  = "Berries"
$ws.Cells.Item(87, 9).Value  = 100101001
$ws.Cells.Item(87, 10).Value = "Arándano (blue)"
$ws.Cells.Item(87, 11).Value = "Sin especificar"
$ws.Cells.Item(87, 12).Value = "Primera"
$ws.Cells.Item(87, 13).Value = 100
$ws.Cells.Item(87, 14).Value = 2800
$ws.Cells.Item(87, 15).Value = 3000
$ws.Cells.Item(87, 16).Value = 2900
$ws.Cells.Item(87, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(87, 18).Value = "Provincia de Colchagua"
$ws.Cells.Item(87, 19).Value = 1450
$ws.Cells.Item(87, 20).Value = 2
